# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# per the latest scrape. D-column writes force text via a temporary
# "@" (Text) number format so Excel does not reinterpret values like
# "311.59" as a floating-point number; the format is reset to Normal
# immediately after so the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.361.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.467.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.72%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.66%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  -4.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.59%  "

$ws.Range("E11").Value = "  -3.03%  "

$ws.Range("E12").Value = "  -1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.847.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.485.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.14%  "

$ws.Range("E17").Value = "  -3.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.352.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0923"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.14%  "

$ws.Range("E21").Value = "  -9.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.89%  "

$ws.Range("E24").Value = "  -4.38%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.99%  "

$ws.Range("E28").Value = "  -4.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "151.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.13%  "

$ws.Range("E32").Value = "  -6.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.56%  "

$ws.Range("E34").Value = "  -3.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0747"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.95%  "

$ws.Range("E37").Value = "  -6.66%  "

$ws.Range("E38").Value = "  -5.76%  "

$ws.Range("E39").Value = "  -3.14%  "

$ws.Range("E40").Value = "  -7.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.14%  "

$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.984.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.10%  "

$ws.Range("E45").Value = "  -4.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.47%  "

$ws.Range("E47").Value = "  -5.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.710.97"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "69.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "96.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.92%  "
